# Edit script: applies the "changes to utility and moments" commit to Outcomes.xlsx
#  - eliminate one moment (row 20 on "data" sheet: "% Intermediate control")
#  - change utility function (updated simulation moments in data!D:F)
#  - added simce measure function (scale corrections *100 on "table" sheet)
#  - fixed aep income (re-derived D/E/F values on "data" sheet)

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")
$tableWs = $wb.Worksheets.Item("table")

# --- 1. Update simulation moment values on the "data" sheet (rows 5-23) ---
$dataWs.Cells.Item(5, 4).Value = 2.439958992473283
$dataWs.Cells.Item(5, 5).Value = 2.473020242929459
$dataWs.Cells.Item(5, 6).Value = 0.07832261502918832
$dataWs.Cells.Item(6, 4).Value = 0.06927355064914854
$dataWs.Cells.Item(6, 5).Value = 0.07376992758363485
$dataWs.Cells.Item(6, 6).Value = 0.002756964555481807
$dataWs.Cells.Item(7, 4).Value = -0.4475215708873302
$dataWs.Cells.Item(7, 5).Value = -0.3552965564727783
$dataWs.Cells.Item(7, 6).Value = 0.01283204546441114
$dataWs.Cells.Item(8, 4).Value = 0.1367171574465806
$dataWs.Cells.Item(8, 5).Value = 0.2189719285666943
$dataWs.Cells.Item(8, 6).Value = 0.008714648232788468
$dataWs.Cells.Item(9, 4).Value = 2.464592548098875
$dataWs.Cells.Item(9, 5).Value = 2.589592592477798
$dataWs.Cells.Item(9, 6).Value = 0.08230323421738787
$dataWs.Cells.Item(10, 4).Value = 0.2942513028456653
$dataWs.Cells.Item(10, 5).Value = 0.3511381005048752
$dataWs.Cells.Item(10, 6).Value = 0.01203229253374716
$dataWs.Cells.Item(11, 4).Value = 2.418578590891703
$dataWs.Cells.Item(11, 5).Value = 2.492724134556088
$dataWs.Cells.Item(11, 6).Value = 0.07914624579444615
$dataWs.Cells.Item(12, 4).Value = 0.02549858700040371
$dataWs.Cells.Item(12, 5).Value = 0.03895805112074458
$dataWs.Cells.Item(12, 6).Value = 0.004028269770326015
$dataWs.Cells.Item(13, 4).Value = 0.3776907549454986
$dataWs.Cells.Item(13, 5).Value = 0.2825256499751223
$dataWs.Cells.Item(13, 6).Value = 0.01278570880286808
$dataWs.Cells.Item(14, 4).Value = 0.3479127977392007
$dataWs.Cells.Item(14, 5).Value = 0.5282867442094009
$dataWs.Cells.Item(14, 6).Value = 0.01946973085298174
$dataWs.Cells.Item(15, 4).Value = 0.2288978603148971
$dataWs.Cells.Item(15, 5).Value = 0.1492295546947322
$dataWs.Cells.Item(15, 6).Value = 0.008675134419246397
$dataWs.Cells.Item(16, 4).Value = 0.1569736613206976
$dataWs.Cells.Item(16, 5).Value = 0.2381830824043836
$dataWs.Cells.Item(16, 6).Value = 0.01986858011483212
$dataWs.Cells.Item(17, 4).Value = 0.1666339721265211
$dataWs.Cells.Item(17, 5).Value = 0.2232180179692148
$dataWs.Cells.Item(17, 6).Value = 0.02229085519208141
$dataWs.Cells.Item(18, 4).Value = -0.0002070351378107086
$dataWs.Cells.Item(18, 5).Value = -0.06485860670546678
$dataWs.Cells.Item(18, 6).Value = 0.01981038070346321
$dataWs.Cells.Item(19, 4).Value = -0.002116344328242216
$dataWs.Cells.Item(19, 5).Value = -0.02306611005363123
$dataWs.Cells.Item(19, 6).Value = 0.02042360966074426
$dataWs.Cells.Item(21, 4).Value = 0.4730797615668464
$dataWs.Cells.Item(21, 5).Value = 0.5851949692469636
$dataWs.Cells.Item(21, 6).Value = 0.02042369864481303
$dataWs.Cells.Item(22, 4).Value = 0.4943662863465367
$dataWs.Cells.Item(22, 5).Value = 0.4185176715820618
$dataWs.Cells.Item(22, 6).Value = 0.01967542804872869
$dataWs.Cells.Item(23, 4).Value = 0.247284977019973
$dataWs.Cells.Item(23, 5).Value = 0.2080090161961711
$dataWs.Cells.Item(23, 6).Value = 0.01423183394444534

# --- 2. Eliminate the "% Intermediate control" moment (row 20): clear its
#        numeric columns but keep the I20 style cell, and mark it ELIMINADO ---
$dataWs.Range("D20:J20").ClearContents()
$dataWs.Range("K20").Value = "ELIMINADO"

# --- 3. "table" sheet: delete the now-eliminated "% intermediate (control
#        group)" row (old row 19). This shifts rows 20-24 up to 19-23 and
#        keeps every data! cross-sheet formula reference intact. ---
$tableWs.Rows("19").Delete()

# --- 4. Rescale the moments that now need a *100 factor ---
$tableWs.Range("F12").Formula = "=data!E13*100"
$tableWs.Range("F13").Formula = "=data!E14*100"
$tableWs.Range("F14").Formula = "=data!E15*100"
$tableWs.Range("F19").Formula = "=data!E21*100"

# --- 5. Selections / active sheet: "table" becomes the active tab, with
#        the selection at H14; "data" selection moves to D21:E21. ---
$dataWs.Range("D21:E21").Select()
$tableWs.Activate()
$tableWs.Range("H14").Select()
